$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Control groups")
$shp = $ws.Shapes.Item(1)
$tr = $shp.TextFrame2.TextRange
$tr.Text = "Hello`r`nWorld`r`nFoo"
